# Release V 1.2.1 Ejemplos de Sintesis terminados
# Rebuild the "Senal Inventada" data + chart: the sampled signal now has
# 34 points (was 32), covering a narrower x-range (1E-4 .. 3.4E-3) and a
# reshaped y-series (rectangular pulse + triangular pulse + trapezoid).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rewrite the source data in columns A (x) and B (y), rows 1-34.
# ---------------------------------------------------------------------
$yValues = @(0, 0.2, 0.4, 0.6, 0.8, 1, 1.2, 1, 0.8, 0.6, 0.4, 0.2, 0, 0.19, 0.3, 0.37, 0.41, 0.43, 0.43, 0.41, 0.37, 0.3, 0.19, 0, 2.7, 2.9, 3.03, 3.1, 3.12, 3.1, 3.03, 2.9, 2.7, 0)

for ($r = 1; $r -le 34; $r++) {
    $ws.Cells.Item($r, 1).Value = $r * 0.0001
    $ws.Cells.Item($r, 2).Value = $yValues[$r - 1]
}

# The old B14 cell carried a one-off "#,##0" number format (style index 1);
# the refreshed sheet no longer singles that cell out, so drop back to the
# plain/default style now that its value is a normal fraction (0.19).
$ws.Range("B14").Style = "Normal"

# ---------------------------------------------------------------------
# 2) Point the chart series at the new, larger A1:B34 range.
# ---------------------------------------------------------------------
$co = $ws.ChartObjects(1)
$chart = $co.Chart

try {
    $chart.SetSourceData($ws.Range("A1:B34"))
} catch {
}

try {
    $ser = $chart.SeriesCollection(1)
    $ser.XValues = $ws.Range("A1:A34")
    $ser.Values = $ws.Range("B1:B34")
    $ser.Formula = "=SERIES(,Hoja1!`$A`$1:`$A`$34,Hoja1!`$B`$1:`$B`$34,1)"
} catch {
}

# ---------------------------------------------------------------------
# 3) Give the plot area an explicit (manual) inner layout, matching the
#    tighter margins the chart ended up with after the resize below.
# ---------------------------------------------------------------------
try {
    $pa = $chart.PlotArea
    $pa.InsideLeft = 0.067330927384076991
    $pa.InsideTop = 0.13930555555555557
    $pa.InsideWidth = 0.8680719597550306
    $pa.InsideHeight = 0.72088764946048411
} catch {
}

# ---------------------------------------------------------------------
# 4) Move + resize the chart (it shrinks a little and shifts up-left,
#    now spanning C6 .. H21 instead of C6 .. I20).
# ---------------------------------------------------------------------
$co.Left = 127.68334645669292
$co.Top = 87.95677165354331
$co.Width = 323.9008858267717
$co.Height = 225.47795275590548
